$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Delete row 16 first (lone "end screen" row), then row 13 (lone "begin screen" row)
# so row indices for the first deletion remain valid.
$ws.Rows.Item(16).Delete()
$ws.Rows.Item(13).Delete()
